$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cell contents: A1 "First Name" -> "FirstName[0]", A2 "Last Name" -> "LastName[0]"
$ws.Range("A1").Value = "FirstName[0]"
$ws.Range("A2").Value = "LastName[0]"

# Update the selected cell to A2 (matches the sheetView selection in the target file)
$ws.Range("A2").Select()
